# Rename the inline picture shapes (Pearson logo x2 in the footers, BTEC
# logo in the first-page header) by swapping their image1.x <-> image2.x
# display names. `InlineShape` has no settable `.Name` in the Word object
# model, so the shape is briefly converted to a floating shape (which does
# expose `.Name`), renamed, then converted back to an inline shape so the
# drawing stays anchored in the text flow exactly as it was.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShape($inlineShape, $newName) {
    $floating = $inlineShape.ConvertToShape()
    $floating.Name = $newName
    $floating.ConvertToInlineShape() | Out-Null
}

# First-page footer (footer1.xml): PearsonLogo image2.png -> image1.png
$footerFirst = $sec.Footers.Item(2)
Rename-InlineShape $footerFirst.Range.InlineShapes.Item(1) "image1.png"

# Default footer (footer2.xml): PearsonLogo image2.png -> image1.png
$footerDefault = $sec.Footers.Item(1)
Rename-InlineShape $footerDefault.Range.InlineShapes.Item(1) "image1.png"

# First-page header (header1.xml): BTec_Logo-Orange image1.jpg -> image2.jpg
$headerFirst = $sec.Headers.Item(2)
Rename-InlineShape $headerFirst.Range.InlineShapes.Item(1) "image2.jpg"
